$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = -0.3422723562191532
$ws.Range("C3").Value = -0.5037688924316441
$ws.Range("E3").Value = -0.467076459743887
$ws.Range("C4").Value = 0.169534172659791
$ws.Range("E4").Value = -0.03968684591929561
$ws.Range("C5").Value = 1.032338390744236
$ws.Range("E5").Value = 0.3765075513336269
$ws.Range("C6").Value = -0.1937612543835177
$ws.Range("E6").Value = 0.6066710853121382
$ws.Range("C7").Value = -0.8331679621937482
$ws.Range("E7").Value = -0.3970496740026364
$ws.Range("C8").Value = 0.1273541662098365
$ws.Range("E8").Value = -0.4865818826308876
$ws.Range("C9").Value = -0.08273351073040391
$ws.Range("E9").Value = -0.01252079199893785
$ws.Range("C10").Value = -0.1151820594382569
$ws.Range("E10").Value = 0.07348980370169844
$ws.Range("C11").Value = 0.157394256377752
$ws.Range("E11").Value = 0.06194937150048041
$ws.Range("C12").Value = 0.3426151435189873
$ws.Range("E12").Value = 0.1686730364466316
$ws.Range("C13").Value = 0.1598952850611068
$ws.Range("E13").Value = 0.132465972367557
$ws.Range("C14").Value = -0.4923796969465988
$ws.Range("E14").Value = -0.3459257698102514
$ws.Range("C15").Value = -0.5121403324772844
$ws.Range("E15").Value = -1.073589070820447
$ws.Range("C16").Value = -0.2454721753057276
$ws.Range("E16").Value = -1.238905350026021
$ws.Range("C17").Value = 0.7038634017465073
$ws.Range("E17").Value = 0.4406734233171727
$ws.Range("C18").Value = 0.20168190406884
$ws.Range("E18").Value = 0.3613321345859122
$ws.Range("C19").Value = -0.06418790329880686
$ws.Range("E19").Value = -0.08988642825158433
